$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 135102.67
$ws.Range("J17").Value = 135102.67
$ws.Range("L17").Value = 405308.01
$ws.Range("N17").Value = -405644.01
$ws.Range("H33").Value = 300.8
$ws.Range("I33").Value = 278.66666
$ws.Range("J33").Value = 500
$ws.Range("K33").Value = 278.66666
$ws.Range("L33").Value = 500
$ws.Range("M33").Value = -49.66665999999998
$ws.Range("N33").Value = -958
$ws.Range("H51").Value = 9469.6
$ws.Range("J51").Value = 10012.667
$ws.Range("L51").Value = 10012.667
$ws.Range("N51").Value = -10980.667
$ws.Range("H96").Value = 2275.25
$ws.Range("I96").Value = 1447.75
$ws.Range("J96").Value = 3102.75
$ws.Range("K96").Value = 4343.25
$ws.Range("L96").Value = 9308.25
$ws.Range("M96").Value = -2970.25
$ws.Range("N96").Value = -12054.25
$ws.Range("H100").Value = 1095.8889
$ws.Range("I100").Value = 1290
$ws.Range("J100").Value = 998.8333
$ws.Range("K100").Value = 1290
$ws.Range("L100").Value = 998.8333
$ws.Range("M100").Value = -749
$ws.Range("N100").Value = -2080.8333
$ws.Range("H116").Value = 7918.222
$ws.Range("J116").Value = 4249.9
$ws.Range("L116").Value = 4249.9
$ws.Range("N116").Value = -11133.9
$ws.Range("H132").Value = 2146.1482
$ws.Range("I132").Value = 2016.4762
$ws.Range("K132").Value = 6049.4286
$ws.Range("M132").Value = -3519.4286
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3386.4412
$ws.Range("I32").Value = 2246.3447
$ws.Range("K32").Value = 2246.3447
$ws.Range("M32").Value = -1959.3447
$ws.Range("H45").Value = 836777.75
$ws.Range("I45").Value = 1113204
$ws.Range("J45").Value = 7499
$ws.Range("K45").Value = 1113204
$ws.Range("L45").Value = 7499
$ws.Range("M45").Value = -1112827
$ws.Range("N45").Value = -8253
$ws.Range("H63").Value = 4047.7
$ws.Range("I63").Value = 3809.75
$ws.Range("K63").Value = 3809.75
$ws.Range("M63").Value = -3123.75
$ws.Range("H66").Value = 4047.7
$ws.Range("I66").Value = 3809.75
$ws.Range("K66").Value = 19048.75
$ws.Range("M66").Value = -15616.75
$ws.Range("H96").Value = 55085.75
$ws.Range("J96").Value = 55085.75
$ws.Range("L96").Value = 55085.75
$ws.Range("N96").Value = -60577.75
$ws.Range("H97").Value = 603.61536
$ws.Range("I97").Value = 568
$ws.Range("J97").Value = 799.5
$ws.Range("K97").Value = 568
$ws.Range("L97").Value = 799.5
$ws.Range("M97").Value = -72
$ws.Range("N97").Value = -1791.5
$ws.Range("H132").Value = 2086652.1
$ws.Range("I132").Value = 2859863
$ws.Range("J132").Value = 4930.3076
$ws.Range("K132").Value = 8579589
$ws.Range("L132").Value = 14790.9228
$ws.Range("M132").Value = -8577059
$ws.Range("N132").Value = -19850.9228
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1944.4828
$ws.Range("I20").Value = 1616.4445
$ws.Range("K20").Value = 1616.4445
$ws.Range("M20").Value = -1369.4445
$ws.Range("H86").Value = 1750.4375
$ws.Range("I86").Value = 1898
$ws.Range("J86").Value = 1425.8
$ws.Range("K86").Value = 1898
$ws.Range("L86").Value = 1425.8
$ws.Range("M86").Value = -775
$ws.Range("N86").Value = -3671.8
$ws.Range("H89").Value = 1750.4375
$ws.Range("I89").Value = 1898
$ws.Range("J89").Value = 1425.8
$ws.Range("K89").Value = 9490
$ws.Range("L89").Value = 7129
$ws.Range("M89").Value = -3874
$ws.Range("N89").Value = -18361
$ws.Range("H100").Value = 21124.75
$ws.Range("J100").Value = 25499.666
$ws.Range("L100").Value = 25499.666
$ws.Range("N100").Value = -27663.666
$ws.Range("H132").Value = 107500
$ws.Range("J132").Value = 107500
$ws.Range("L132").Value = 107500
$ws.Range("N132").Value = -117620
$ws.Range("H134").Value = 51003204
$ws.Range("I134").Value = 85001930
$ws.Range("K134").Value = 255005790
$ws.Range("M134").Value = -255003255
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 217.6923
$ws.Range("I7").Value = 64.57143000000001
$ws.Range("K7").Value = 64.57143000000001
$ws.Range("M7").Value = 48.42856999999999
$ws.Range("H22").Value = 8600.416999999999
$ws.Range("I22").Value = 14510.857
$ws.Range("J22").Value = 325.8
$ws.Range("K22").Value = 14510.857
$ws.Range("L22").Value = 325.8
$ws.Range("M22").Value = -14160.857
$ws.Range("N22").Value = -1025.8
$ws.Range("H31").Value = 10092.172
$ws.Range("I31").Value = 6204.3887
$ws.Range("K31").Value = 6204.3887
$ws.Range("M31").Value = -5909.3887
$ws.Range("H34").Value = 10092.172
$ws.Range("I34").Value = 6204.3887
$ws.Range("K34").Value = 6204.3887
$ws.Range("M34").Value = -6002.3887
$ws.Range("H105").Value = 2223842.8
$ws.Range("J105").Value = 2440.3333
$ws.Range("L105").Value = 2440.3333
$ws.Range("N105").Value = -5934.3333
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 441.16666
$ws.Range("I2").Value = 597.25
$ws.Range("J2").Value = 129
$ws.Range("K2").Value = 3583.5
$ws.Range("L2").Value = 774
$ws.Range("M2").Value = -3470.5
$ws.Range("N2").Value = -1000
$ws.Range("H17").Value = 586
$ws.Range("J17").Value = 464
$ws.Range("L17").Value = 1392
$ws.Range("N17").Value = -1730
$ws.Range("H34").Value = 1158.0555
$ws.Range("J34").Value = 1984.1111
$ws.Range("L34").Value = 5952.3333
$ws.Range("N34").Value = -6120.3333
$ws.Range("H39").Value = 2124.5
$ws.Range("J39").Value = 3250
$ws.Range("L39").Value = 9750
$ws.Range("N39").Value = -10338
$ws.Range("H55").Value = 2650
$ws.Range("J55").Value = 5000
$ws.Range("L55").Value = 15000
$ws.Range("N55").Value = -15354
$ws.Range("H87").Value = 14416.667
$ws.Range("I87").Value = 14416.667
$ws.Range("K87").Value = 43250.001
$ws.Range("M87").Value = -42002.001
$ws.Range("H90").Value = 14416.667
$ws.Range("I90").Value = 14416.667
$ws.Range("K90").Value = 129750.003
$ws.Range("M90").Value = -123510.003
$ws.Range("H107").Value = 1361.1666
$ws.Range("J107").Value = 1848.8572
$ws.Range("L107").Value = 5546.571599999999
$ws.Range("N107").Value = -9386.571599999999
$ws.Range("H121").Value = 67967.734
$ws.Range("I121").Value = 125371
$ws.Range("K121").Value = 376113
$ws.Range("M121").Value = -374803
$ws.Range("H122").Value = 831.6667
$ws.Range("I122").Value = 195.875
$ws.Range("K122").Value = 1762.875
$ws.Range("M122").Value = 687.125
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 168.10527
$ws.Range("I2").Value = 77.454544
$ws.Range("J2").Value = 292.75
$ws.Range("K2").Value = 77.454544
$ws.Range("L2").Value = 292.75
$ws.Range("M2").Value = 35.545456
$ws.Range("N2").Value = -518.75
$ws.Range("H80").Value = 3598.4443
$ws.Range("I80").Value = 3672.5
$ws.Range("J80").Value = 3539.2
$ws.Range("K80").Value = 3672.5
$ws.Range("L80").Value = 3539.2
$ws.Range("M80").Value = -2674.5
$ws.Range("N80").Value = -5535.2
$ws.Range("H83").Value = 3598.4443
$ws.Range("I83").Value = 3672.5
$ws.Range("J83").Value = 3539.2
$ws.Range("K83").Value = 18362.5
$ws.Range("L83").Value = 17696
$ws.Range("M83").Value = -13370.5
$ws.Range("N83").Value = -27680
$ws.Range("H102").Value = 1135
$ws.Range("I102").Value = 1068.4615
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 1068.4615
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = 553.5385000000001
$ws.Range("N102").Value = -5244
$ws.Range("H132").Value = 6251772
$ws.Range("I132").Value = 6946324.5
$ws.Range("J132").Value = 797.5
$ws.Range("K132").Value = 20838973.5
$ws.Range("L132").Value = 2392.5
$ws.Range("M132").Value = -20836443.5
$ws.Range("N132").Value = -7452.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 541.1667
$ws.Range("I55").Value = 500
$ws.Range("J55").Value = 549.4
$ws.Range("K55").Value = 500
$ws.Range("L55").Value = 549.4
$ws.Range("M55").Value = -327
$ws.Range("N55").Value = -895.4
$ws.Range("H136").Value = 1777.3214
$ws.Range("I136").Value = 1494.1666
$ws.Range("K136").Value = 4482.4998
$ws.Range("M136").Value = -1932.4998
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 15831.667
$ws.Range("I74").Value = 15998.5
$ws.Range("J74").Value = 15748.25
$ws.Range("K74").Value = 15998.5
$ws.Range("L74").Value = 15748.25
$ws.Range("M74").Value = -15062.5
$ws.Range("N74").Value = -17620.25
$ws.Range("H77").Value = 15831.667
$ws.Range("I77").Value = 15998.5
$ws.Range("J77").Value = 15748.25
$ws.Range("K77").Value = 47995.5
$ws.Range("L77").Value = 47244.75
$ws.Range("M77").Value = -43315.5
$ws.Range("N77").Value = -56604.75
$ws.Range("H96").Value = 2723.875
$ws.Range("I96").Value = 3586.625
$ws.Range("J96").Value = 1861.125
$ws.Range("K96").Value = 3586.625
$ws.Range("L96").Value = 1861.125
$ws.Range("M96").Value = -2213.625
$ws.Range("N96").Value = -4607.125
$ws.Range("H126").Value = 1099.1072
$ws.Range("I126").Value = 1045.04
$ws.Range("K126").Value = 3135.12
$ws.Range("M126").Value = -665.1199999999999
$ws.Range("H132").Value = 21741828
$ws.Range("I132").Value = 33335904
$ws.Range("J132").Value = 2938
$ws.Range("K132").Value = 100007712
$ws.Range("L132").Value = 8814
$ws.Range("M132").Value = -100005182
$ws.Range("N132").Value = -13874
